$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Insert two new rows after the current row 4 (new rows 5 and 6)
#    This pushes the old row 6 (italic note) and row 7 (final note)
#    down to row 8 and row 9 respectively.
# ---------------------------------------------------------------
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# ---------------------------------------------------------------
# 2. Header row (row 1) - add the new "Person" column header in F1
# ---------------------------------------------------------------
$ws.Cells.Item(1,6).Value = "Person"

# ---------------------------------------------------------------
# 3. Row 2 stays S102 / Semakau Landfill, but remark text + new Person column
# ---------------------------------------------------------------
$ws.Cells.Item(2,5).Value = "No HDB town can be mapped; Included in station table but excluded from locations table"
$ws.Cells.Item(2,6).Value = "Erica"

# ---------------------------------------------------------------
# 4. Row 3 becomes S115 / Tuas South Avenue 3 (previously row 4's data)
# ---------------------------------------------------------------
$ws.Cells.Item(3,3).Value = "S115"
$ws.Cells.Item(3,4).Value = "Tuas South Avenue 3"
$ws.Cells.Item(3,5).Value = "No HDB town can be mapped; Included in station table but excluded from locations table"
$ws.Cells.Item(3,6).Value = "Erica"

# ---------------------------------------------------------------
# 5. Row 4 becomes S106 / Pulau Ubin (previously row 3's data), with a
#    change date of 45622 instead of 45621
# ---------------------------------------------------------------
$ws.Cells.Item(4,1).Value = 45622
$ws.Cells.Item(4,3).Value = "S106"
$ws.Cells.Item(4,4).Value = "Pulau Ubin"
$ws.Cells.Item(4,5).Value = "No data on air_temp table and humidity table; Included in station table but excluded from locations table"
$ws.Cells.Item(4,6).Value = "Erica"

# ---------------------------------------------------------------
# 6. New row 5 - floor_area_sqm change
# ---------------------------------------------------------------
$ws.Cells.Item(5,1).Value = 45622
$ws.Cells.Item(5,2).Value = "floor_area_sqm"
$ws.Cells.Item(5,3).Value = "(n/a)"
$ws.Cells.Item(5,5).Value = "Change data type to float, as having ingestion issue via Python"
$ws.Cells.Item(5,6).Value = "Yvonne"

# ---------------------------------------------------------------
# 7. New row 6 - resale_price change
# ---------------------------------------------------------------
$ws.Cells.Item(6,1).Value = 45622
$ws.Cells.Item(6,2).Value = "resale_price"
$ws.Cells.Item(6,3).Value = "(n/a)"
$ws.Cells.Item(6,5).Value = "Change data type to float, as having ingestion issue via Python"
$ws.Cells.Item(6,6).Value = "Yvonne"

# Make sure date formatting (m/d/yyyy, shared with A2:A4) is applied on A5:A6
$ws.Cells.Item(4,1).Copy()
$ws.Range($ws.Cells.Item(5,1), $ws.Cells.Item(6,1)).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# 8. Row 8 (previously row 6) - the italic explanation note now also
#    has a backtick placed in column E
# ---------------------------------------------------------------
$ws.Cells.Item(8,5).Value = "``"

# Row 9 (previously row 7) keeps its original text - nothing to change

# ---------------------------------------------------------------
# 9. Column widths: E narrower, new column F added
#    (COM ColumnWidth is quantized to the nearest pixel internally, so
#    use the input values that land closest to the target stored widths)
# ---------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 82.65
$ws.Columns.Item(6).ColumnWidth = 11.25

# ---------------------------------------------------------------
# 10. Move the picture down by 2 rows (28.8 points at 14.4pt/row) to
#     keep it anchored below the table, matching the row insert shift
# ---------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Top = $shp.Top + 28.8

# ---------------------------------------------------------------
# 11. Reselect cell A1 so the sheet view no longer references the old
#     B3 selection
# ---------------------------------------------------------------
$ws.Cells.Item(1,1).Select()
